# The commit replaces the page-count mentioned in the intro paragraph
# ("... terdapat 5 halaman." -> "... terdapat 6 halaman.") The rest of the
# underlying XML diff (removed w:proofErr spell-check markers and merged
# runs) carries no visible text change and is simply how Word
# re-serializes the paragraph after it is edited/re-saved, so a single
# targeted Find & Replace on the visible text reproduces the effective
# edit.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "terdapat 5 halaman.",   # FindText
    $true,                   # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "terdapat 6 halaman.",   # ReplaceWith
    2                        # Replace (wdReplaceAll)
)
